# Työaikakirjanpito update: add a new work entry row, extend the table/sheet
# to 60 rows, and move the running-total SUM formula down to row 60.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Clear the old "Total" / SUM / "Done!" row (row 39). Clearing the
#    contents keeps the existing per-cell styles (date fmt / number fmt /
#    centered general) which already match what the new blank rows need.
# ---------------------------------------------------------------------
$ws.Range("A39:C39").ClearContents()

# ---------------------------------------------------------------------
# 2. Turn the old blank spacer row (row 38) into a new data row for the
#    latest work session.
# ---------------------------------------------------------------------
$ws.Cells.Item(38, 1).Value = 44063
$ws.Cells.Item(38, 2).Value = 6
$ws.Cells.Item(38, 3).Value = "Change password ja reset password toimintaan sekä frontin layout suunnittelua."
$ws.Rows.Item(38).RowHeight = 30

# ---------------------------------------------------------------------
# 3. Add new blank rows 40-59, matching the formatting used by the other
#    blank rows in the sheet (date format / accounting number format /
#    centered general, respectively for columns A/B/C).
# ---------------------------------------------------------------------
for ($r = 40; $r -le 59; $r++) {
    $a = $ws.Cells.Item($r, 1)
    $a.NumberFormat = "[$-409]d\-mmm;@"
    $a.HorizontalAlignment = -4108
    $a.VerticalAlignment = -4108

    $b = $ws.Cells.Item($r, 2)
    $b.NumberFormat = "0.00;[Red]0.00"
    $b.HorizontalAlignment = -4108
    $b.VerticalAlignment = -4108

    $c = $ws.Cells.Item($r, 3)
    $c.HorizontalAlignment = -4108
    $c.VerticalAlignment = -4108
}

# ---------------------------------------------------------------------
# 4. Add the new running-total row 60 with the SUM formula moved down.
# ---------------------------------------------------------------------
$a60 = $ws.Cells.Item(60, 1)
$a60.NumberFormat = "[$-409]d\-mmm;@"
$a60.HorizontalAlignment = -4108
$a60.VerticalAlignment = -4108

$b60 = $ws.Cells.Item(60, 2)
$b60.Formula = "=SUM(B`$2:B`$59)"
$b60.NumberFormat = "0.00;[Red]0.00"
$b60.HorizontalAlignment = -4108
$b60.VerticalAlignment = -4108

$c60 = $ws.Cells.Item(60, 3)
$c60.HorizontalAlignment = -4108
$c60.VerticalAlignment = -4108

# ---------------------------------------------------------------------
# 5. Resize the table / autofilter to cover the new range.
# ---------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:C60"))

# ---------------------------------------------------------------------
# 6. Update the view: scroll down a bit and select B61 (the cell right
#    below the new total row), matching where the user was working.
# ---------------------------------------------------------------------
$ws.Activate() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 31
$win.ScrollColumn = 1
$ws.Range("B61").Select() | Out-Null

Write-Output "done"
